$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.23538613319397
$ws.Range("B1").Value = 2.626508712768555
$ws.Range("C1").Value = 5.109991550445557
$ws.Range("D1").Value = 2.032225370407104
$ws.Range("E1").Value = 1.174598217010498
